$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(1.0342080460121794, 0.7167363786836042, 1.7073445035147568, 3.0450042025624375)"
$ws.Range("C2").Value = "NIG(1.3494066382761636, 1.0222525412030188, 4.0077496499473595, 5.981840650360904)"
$ws.Range("D2").Value = "JSU(-1.171071876238241, 1.0503248436437749, -0.2562736652361195, 2.347183198132743)"
$ws.Range("E2").Value = "NIG(1.4224486559263763, 1.0885506991715033, 3.4749120890520224, 5.265986953654944)"
